$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'266.77"
$ws.Range("D3").Value = "'21.38"
$ws.Range("D4").Value = "'6.119"
$ws.Range("D5").Value = "'0.06102"
$ws.Range("D6").Value = "'3.574"
$ws.Range("D7").Value = "'6.487"
$ws.Range("D9").Value = "'0.8207"
$ws.Range("D10").Value = "'0.01344"
$ws.Range("D11").Value = "'0.1587"
$ws.Range("D12").Value = "'0.08050"
$ws.Range("D13").Value = "'0.03444"
$ws.Range("D15").Value = "'0.09220"
$ws.Range("D16").Value = "'3.759"
$ws.Range("D17").Value = "'0.001629"
$ws.Range("D18").Value = "'0.04644"
$ws.Range("D19").Value = "'0.006390"
$ws.Range("D20").Value = "'0.006146"
$ws.Range("D21").Value = "'0.001069"
$ws.Range("D24").Value = "'2.296"
$ws.Range("D25").Value = "'0.3312"
$ws.Range("D40").Value = "'0.04588"
$ws.Range("D41").Value = "'0.006994"
$ws.Range("D43").Value = "'0.1116"
$ws.Range("D44").Value = "'0.01175"
$ws.Range("D45").Value = "'0.00005846"
$ws.Range("D46").Value = "'0.0009902"
$ws.Range("D49").Value = "'0.001125"
